$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Category" column (D) has some blank cells which break the sheet's
# filters. Fill every row that has data (column C / youtuber name) but is
# missing a Category value with the fallback category "Mixed".
$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $categoryCell = $ws.Cells.Item($r, 3)
    $mixedCell = $ws.Cells.Item($r, 4)
    if (($null -ne $categoryCell.Value2) -and ($null -eq $mixedCell.Value2)) {
        $mixedCell.Value2 = "Mixed"
    }
}

# Reflect the cell the author was last looking at when saving the file.
$ws.Range("D8").Select()
